$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto table figures. A handful of Price cells in column D
# are plain decimal-looking strings (e.g. '262.41'); Excel would normally
# auto-convert those to numbers on assignment, so we briefly force the
# cell to Text format, write the literal string, then restore the original
# 'General' number format so the cell's look-and-feel is unchanged.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.496.05'
$ws.Range("D2").NumberFormat = "General"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.017.51'
$ws.Range("D3").NumberFormat = "General"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.41'
$ws.Range("D5").NumberFormat = "General"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.57%  '
$ws.Range("E5").NumberFormat = "General"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E6").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.97'
$ws.Range("D8").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -7.36%  '
$ws.Range("E8").NumberFormat = "General"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("E9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0773'
$ws.Range("D10").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("E10").NumberFormat = "General"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("E11").NumberFormat = "General"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.92%  '
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.312.53'
$ws.Range("D13").NumberFormat = "General"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("E13").NumberFormat = "General"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.14%  '
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.95'
$ws.Range("D15").NumberFormat = "General"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -8.23%  '
$ws.Range("E15").NumberFormat = "General"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.90%  '
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.026.54'
$ws.Range("D17").NumberFormat = "General"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.390.44'
$ws.Range("D18").NumberFormat = "General"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.72'
$ws.Range("D19").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("E19").NumberFormat = "General"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.19'
$ws.Range("D21").NumberFormat = "General"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.50'
$ws.Range("D22").NumberFormat = "General"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.08%  '
$ws.Range("E23").NumberFormat = "General"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E24").NumberFormat = "General"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E25").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.02'
$ws.Range("D26").NumberFormat = "General"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("E26").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("D27").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.51%  '
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.69'
$ws.Range("D28").NumberFormat = "General"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E28").NumberFormat = "General"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -11.51%  '
$ws.Range("E29").NumberFormat = "General"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("D30").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("E30").NumberFormat = "General"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").NumberFormat = "General"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.42%  '
$ws.Range("E32").NumberFormat = "General"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0649'
$ws.Range("D33").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.88%  '
$ws.Range("E33").NumberFormat = "General"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E34").NumberFormat = "General"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.82'
$ws.Range("D36").NumberFormat = "General"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("E36").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("D38").NumberFormat = "General"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("E38").NumberFormat = "General"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.18'
$ws.Range("D39").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.64%  '
$ws.Range("E39").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.05'
$ws.Range("D40").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.23'
$ws.Range("D41").NumberFormat = "General"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.81%  '
$ws.Range("E41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0942'
$ws.Range("D42").NumberFormat = "General"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.19%  '
$ws.Range("E42").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.398.53'
$ws.Range("D44").NumberFormat = "General"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.55%  '
$ws.Range("E44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.35'
$ws.Range("D45").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.71'
$ws.Range("D46").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.92%  '
$ws.Range("E46").NumberFormat = "General"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.06'
$ws.Range("D48").NumberFormat = "General"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.38%  '
$ws.Range("E48").NumberFormat = "General"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("E49").NumberFormat = "General"

$ws.Range("B50").Value = 'RocketPoolETH'

$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.203.94'
$ws.Range("D50").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("E50").NumberFormat = "General"

$ws.Range("B51").Value = 'NEARProtocol'

$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.97'
$ws.Range("D51").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.86%  '
$ws.Range("E51").NumberFormat = "General"
